# Applies the "nuevos experimentos no convexos" edit: refreshed numeric
# results for a non-convex bilevel-programming experiment.
#
# In the source workbook every "numeric-looking" value on sheets 1-6 is
# actually stored as TEXT (shared string), only Vector_Alpha!A2:A3 holds
# real numbers. To reproduce that faithfully through the Excel object
# model we flip the target cell to Text format ("@") before writing the
# string, then restore the cell's style back to the sheet's default so no
# visible formatting change remains - only the underlying stored type
# (text vs number) differs, matching the original file.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $ws.Range("A1").Style
}

# NOTE: sheet names "Vector_bf" and "Vector_BF" differ only by case, and
# Worksheets.Item(<name>) resolves case-insensitively, so by-name lookup
# would silently grab the wrong sheet. Index into Worksheets by position
# (1-based, matching workbook.xml <sheets> order) instead.
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Restricciones_del_follower ---------------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3 "A2" "-12.217254528122025 - 2x_1 + 1.128693994280267y_1 - 0.1763584366062918y_2"
Set-TextValue $ws3 "B2" "14.717254528122025"
Set-TextValue $ws3 "D2" "0.92"
Set-TextValue $ws3 "E2" "0"

Set-TextValue $ws3 "A3" "18.33341277407054 + x_1 - 3x_2 - 0.15252621544327932y_1 + 0.023832221163012424y_2"
Set-TextValue $ws3 "B3" "-20.33341277407054"
Set-TextValue $ws3 "D3" "0.36"
Set-TextValue $ws3 "E3" "0"
Set-TextValue $ws3 "F3" "0"

Set-TextValue $ws3 "A4" "99.73384175405148 - 0.9761677788369876y_1 + 0.15252621544327932y_2"
Set-TextValue $ws3 "B4" "-99.73384175405148"
Set-TextValue $ws3 "D4" "0.49"
Set-TextValue $ws3 "E4" "0"

Set-TextValue $ws3 "A5" "-15.753412774070544 + 0.15252621544327932y_1 - 0.023832221163012424y_2"
Set-TextValue $ws3 "B5" "15.583412774070544"
Set-TextValue $ws3 "D5" "0.75"
Set-TextValue $ws3 "E5" "0"
Set-TextValue $ws3 "F5" "0"

# --- Punto_modificado ---------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4 "A2" "51.550000000000004"
Set-TextValue $ws4 "B2" "18.099999999999998"
Set-TextValue $ws4 "C2" "102.2"
Set-TextValue $ws4 "D2" "0.2"

# --- Vector_bf -----------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5 "A2" "3.3804385128693992"
Set-TextValue $ws5 "A3" "-0.9031935176358435"

# --- Vector_BF -----------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6 "A2" "2.0"
Set-TextValue $ws6 "A3" "-1.0"
Set-TextValue $ws6 "A4" "-0.5"
Set-TextValue $ws6 "A5" "-0.0"

# --- Vector_Alpha (true numeric cells) -----------------------------------
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("A2").Value = 0.15000000000000002
$ws7.Range("A3").Value = 0.96
